# Fruta / hortaliza, semanal
#
# A new weekly data row is inserted at row 87 (pushing the existing rows
# 87..167 down to 88..168), and the new row is populated with this week's
# reading for Berenjena at "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 87, shifting everything below it down by one
# (mirrors Excel's Rows("87:87").Insert Shift:=xlShiftDown).
$ws.Rows("87:87").Insert(4)

# Populate the newly inserted row with the new weekly record.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("D87").Value = 45033
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112001
$ws.Range("G87").Value = "Berenjena"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 150
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 9000
$ws.Range("M87").Value = 9000
$ws.Range("N87").Value = "`$/caja 50 unidades"
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 180
$ws.Range("Q87").Value = 50
$ws.Range("R87").Value = "Hortaliza"
